$wb = $excel.ActiveWorkbook

# New "FOR PERIOD: FROM / TO" timestamp (Excel serial date) applied to all sheets.
$newDate = 44305.361805555556

# ---------------------------------------------------------------------------
# Edwaleni
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Edwaleni")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F23:F28").Value = 2.4
$ws.Range("F29").Value = 14.6
$ws.Range("F33:F40").Value = 14.6
$ws.Range("F43:F46").Value = 2.4

# ---------------------------------------------------------------------------
# Maguduza
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Maguduza")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F29").Value = 5.6
$ws.Range("F33:F40").Value = 5.6

# ---------------------------------------------------------------------------
# Ezulwini
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ezulwini")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F23:F28").Value = 20
$ws.Range("F29").Value = 20
$ws.Range("F33:F40").Value = 20
$ws.Range("F43:F46").Value = 20

# ---------------------------------------------------------------------------
# Maguga (F14/H14 were formulas referencing Edwaleni; now plain values)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Maguga")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F29").Value = 10
$ws.Range("F30:F32").Value = 20
$ws.Range("F41:F42").Value = 20
